$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '66.448.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.51%  '
$ws.Range('D3').Value = "'" + '3.253.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.14%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'" + '581.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.35%  '
$ws.Range('D6').Value = "'" + '153.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.80%  '
$ws.Range('D7').Value = "'" + '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'" + '3.246.65'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.15%  '
$ws.Range('E9').Value = '  +5.45%  '
$ws.Range('E10').Value = '  +6.80%  '
$ws.Range('D11').Value = "'" + '0.164'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.66%  '
$ws.Range('D12').Value = "'" + '0.489'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.60%  '
$ws.Range('D13').Value = "'" + '37.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Value = "'" + '0.0000235'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.49%  '
$ws.Range('D15').Value = "'" + '3.784.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.44%  '
$ws.Range('D16').Value = "'" + '66.567.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = "'" + '553.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +12.00%  '
$ws.Range('D18').Value = "'" + '3.258.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.97%  '
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('E20').Value = '  +5.65%  '
$ws.Range('D21').Value = "'" + '14.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.56%  '
$ws.Range('D22').Value = "'" + '0.745'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.03%  '
$ws.Range('D23').Value = "'" + '7.86'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.50%  '
$ws.Range('D24').Value = "'" + '13.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.74%  '
$ws.Range('D25').Value = "'" + '81.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.94%  '
$ws.Range('D26').Value = "'" + '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = "'" + '9.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +18.07%  '
$ws.Range('E28').Value = '  +7.78%  '
$ws.Range('E29').Value = '  +5.96%  '
$ws.Range('D30').Value = "'" + '27.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.09%  '
$ws.Range('D31').Value = "'" + '2.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.94%  '
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').Value = "'" + '1.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.54%  '
$ws.Range('D34').Value = "'" + '566.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.83%  '
$ws.Range('D35').Value = "'" + '5.69'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.23%  '
$ws.Range('D36').Value = "'" + '6.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.72%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = "'" + '55.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.31%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'" + '0.0455'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.14%  '
$ws.Range('D39').Value = "'" + '0.0862'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.02%  '
$ws.Range('D40').Value = "'" + '0.130'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.90%  '
$ws.Range('D41').Value = "'" + '3.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.16%  '
$ws.Range('D42').Value = "'" + '3.206.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.56%  '
$ws.Range('D43').Value = "'" + '8.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('D44').Value = "'" + '0.278'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.37%  '
$ws.Range('D45').Value = "'" + '2.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.49%  '
$ws.Range('D46').Value = "'" + '26.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.53%  '
$ws.Range('D48').Value = "'" + '0.0₃0558'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').Value = "'" + '126.29'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.63%  '
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('D51').Value = "'" + '2.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.64%  '
